# Update "Periodo Mora" column (E16:E23) - previous accounting periods are
# replaced with the new periods, in reverse chronological order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2301"
$ws.Range("E17").Value = "2212"
$ws.Range("E18").Value = "2211"
$ws.Range("E19").Value = "2210"
$ws.Range("E20").Value = "2209"
$ws.Range("E21").Value = "2208"
$ws.Range("E22").Value = "2207"
$ws.Range("E23").Value = "2206"

# Update "Valor Mora" (F) and "Salario Basico" (G) figures for the refreshed
# database rows.
$ws.Range("F16").Value = 33333
$ws.Range("G16").Value = 908526

$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 908526

$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 908526

$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 908526

$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 908526

$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 908526

$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 908526

$ws.Range("F23").Value = 40000
$ws.Range("G23").Value = 908526
